$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 694 entirely; this shifts all rows below it up by one.
$ws.Rows(694).Delete()
